$d = $word.ActiveDocument

# --- 1. Pengisian data siswa (NISN) di tabel Lembar Pengesahan ---
$t = $d.Tables.Item(1)

# Lebarkan kolom NISN (kolom ke-3) dari 878 twips (43.9pt) menjadi 2736 twips (136.8pt)
$t.Columns.Item(3).Width = 136.8

# Isi nomor NISN untuk baris-baris siswa yang sudah diketahui datanya
# Baris 1 = header (Nama/Kelas/NISN)
# Baris 3 = Muhammad Zidan Rozaky
$t.Cell(3, 3).Range.Text = "0066115290"
# Baris 5 = Saeful Anwar
$t.Cell(5, 3).Range.Text = "0063495936"
# Baris 7 = Zidane Satria Perkasa
$t.Cell(7, 3).Range.Text = "0063495936"

# --- 2. Bersihkan paragraf kosong ber-format Heading2 setelah judul "TEMPAT PKL" ---
# (Catatan: memakai $d.Content.Paragraphs, bukan $d.Paragraphs, karena koleksi
#  $d.Paragraphs menjadi tidak dapat diandalkan setelah tabel diakses di atas.)
$paras = $d.Content.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq "TEMPAT PKL`r" -and $p.Style.NameLocal -eq "Heading 2") {
        $nextPara = $paras.Item($i + 1)
        if ($nextPara.Range.Text -eq "`r" -and $nextPara.Style.NameLocal -eq "Heading 2") {
            $nextPara.Style = "Normal"
        }
        break
    }
}
